# Apply updated dSF (column F) values for the listed rows.
# Mapping of row -> new value, per the commit "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    5  = -3
    7  = 3
    8  = -4
    10 = 1
    12 = -7
    13 = 4
    15 = -7
    17 = -4
    18 = -1
    19 = -4
    22 = 0
    24 = -6
    26 = 7
    27 = -5
    30 = -4
    33 = -7
    38 = -4
    40 = -7
    41 = -6
    43 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
